$d = $word.ActiveDocument

# Locate the sentence fragment that needs to be split into three runs:
# " кадров 5 Гц." -> " кадров " + "5" + " Гц."
$rng = $d.Content
$found = $rng.Find.Execute(" кадров 5 Гц.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $fullText = $rng.Text
    $baseStart = $rng.Start

    # Locate "5" inside the matched fragment so we don't rely on hard-coded offsets.
    $numStr = "5"
    $idxNum = $fullText.IndexOf($numStr)

    $numStart = $baseStart + $idxNum
    $numEnd = $numStart + $numStr.Length

    $tailStart = $numEnd
    $tailEnd = $baseStart + $fullText.Length

    # Force Word to materialize run boundaries around "5" and " Гц." by
    # toggling a character property on and back off again, without
    # actually changing the visible formatting (sz/szCs/lang stay as-is).
    $numRange = $d.Range($numStart, $numEnd)
    $numRange.Bold = 1
    $numRange.Bold = 0

    $tailRange = $d.Range($tailStart, $tailEnd)
    $tailRange.Bold = 1
    $tailRange.Bold = 0
}
